$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two shared-string sentences (A2 / B2)
$ws.Range("A2").Value = "não estou usando artigos nessa frase"
$ws.Range("B2").Value = "estou usando artigos nessa frase"

# Move the active-cell selection from B9 to B8
$ws.Range("B8").Select() | Out-Null

# Widen column A (closest achievable value in this runtime's quantized
# ColumnWidth model to the target OOXML width of 36.7109375)
$ws.Columns.Item(1).ColumnWidth = 35.83
